$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '66.683.81'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  -0.24%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.229.64'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  +0.60%  '
$ws.Range("E4").Value = '  -0.10%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '605.41'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -0.41%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '156.14'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -1.36%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.999'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  -0.09%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '3.228.33'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +0.62%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.544'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -2.11%  '
$ws.Range("E10").Value = '  -0.40%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '5.79'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -3.84%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.503'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -2.65%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.0000269'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +0.02%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '38.74'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -1.98%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '3.756.26'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +0.45%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '66.686.06'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -0.23%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '3.227.30'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +0.34%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '7.27'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -3.03%  '
$ws.Range("E19").Value = '  +0.88%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '506.97'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -2.52%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '15.29'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -1.12%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.736'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -1.16%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '8.03'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -2.16%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '14.55'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -4.15%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '85.17'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -0.39%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '0.167'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +83.38%  '
$ws.Range("E27").Value = '  +0.06%  '
$ws.Range("E28").Value = '  -1.50%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '9.03'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -3.66%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '2.37'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -1.22%  '
$ws.Range("B31").Value = 'NEARProtocol'
$ws.Range("C31").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '6.98'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -0.75%  '
$ws.Range("B32").Value = 'Stacks'
$ws.Range("C32").Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '2.91'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -4.20%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '28.36'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -0.15%  '
$ws.Range("E34").Value = '  -0.06%  '
$ws.Range("E35").Value = '  -5.55%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '6.43'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -2.78%  '
$ws.Range("B37").Value = 'Bittensor'
$ws.Range("C37").Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '503.94'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -3.64%  '
$ws.Range("B38").Value = 'OKB'
$ws.Range("C38").Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '55.23'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +0.35%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.0₃0776'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +12.23%  '
$ws.Range("B40").Value = 'Kaspa'
$ws.Range("C40").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.130'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +1.79%  '
$ws.Range("B41").Value = 'VeChain'
$ws.Range("C41").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.0418'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -3.07%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '3.02'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +2.72%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '8.73'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -2.59%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.296'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -2.83%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '2.950.90'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +1.36%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '2.46'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -2.54%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '28.26'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -2.49%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '2.42'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -0.27%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.118'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -0.53%  '
$ws.Range("E50").Value = '  +0.01%  '
$ws.Range("B51").Value = 'CoreDAO'
$ws.Range("C51").Value = 'https://coinranking.com/coin/HFvoXUQh4+coredao-core'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '2.56'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -5.32%  '
